$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous used range so stale cells from the old layout do not linger
$ws.Range("A1:C95").ClearContents()

$ws.Range("A1").Value = '''what'
$ws.Range("B1").Value = '''level'
$ws.Range("C1").Value = '''Overall'
$ws.Range("A2").Value = '''n'
$ws.Range("C2").Value = '''68859'
$ws.Range("A3").Value = '''P_Age (mean (sd))'
$ws.Range("C3").Value = '''68.79 (9.93)'
$ws.Range("A4").Value = '''P_Gender (%)'
$ws.Range("B4").Value = '''Man'
$ws.Range("C4").Value = '''30522 (44.3) '
$ws.Range("B5").Value = '''Kvinna'
$ws.Range("C5").Value = '''38337 (55.7) '
$ws.Range("A6").Value = '''P_BMI (mean (sd))'
$ws.Range("C6").Value = '''27.32 (4.38)'
$ws.Range("A7").Value = '''P_ASA (%)'
$ws.Range("B7").Value = '''1'
$ws.Range("C7").Value = '''16784 (24.4) '
$ws.Range("B8").Value = '''2'
$ws.Range("C8").Value = '''41294 (60.0) '
$ws.Range("B9").Value = '''3'
$ws.Range("C9").Value = '''10781 (15.7) '
$ws.Range("A10").Value = '''P_TypeOfHospital (%)'
$ws.Range("B10").Value = '''Universitets- eller regionssjukhus'
$ws.Range("C10").Value = ''' 4750 ( 6.9) '
$ws.Range("B11").Value = '''Länsdelsjukhus'
$ws.Range("C11").Value = '''27698 (40.2) '
$ws.Range("B12").Value = '''Länssjukhus'
$ws.Range("C12").Value = '''21732 (31.6) '
$ws.Range("B13").Value = '''Privatsjukhus'
$ws.Range("C13").Value = '''14679 (21.3) '
$ws.Range("A14").Value = '''P_ProtGrp (%)'
$ws.Range("B14").Value = '''Cemented'
$ws.Range("C14").Value = '''46115 (67.0) '
$ws.Range("B15").Value = '''Cementless'
$ws.Range("C15").Value = '''11731 (17.0) '
$ws.Range("B16").Value = '''Hybrid'
$ws.Range("C16").Value = ''' 1929 ( 2.8) '
$ws.Range("B17").Value = '''Reversed hybrid'
$ws.Range("C17").Value = ''' 9084 (13.2) '
$ws.Range("A18").Value = '''P_SurgYear (mean (sd))'
$ws.Range("C18").Value = ''' 3.74 (2.26)'
$ws.Range("A19").Value = '''education (%)'
$ws.Range("B19").Value = '''low'
$ws.Range("C19").Value = '''22815 (33.1) '
$ws.Range("B20").Value = '''high'
$ws.Range("C20").Value = '''17671 (25.7) '
$ws.Range("B21").Value = '''middle'
$ws.Range("C21").Value = '''28373 (41.2) '
$ws.Range("A22").Value = '''civil_status (%)'
$ws.Range("B22").Value = '''married'
$ws.Range("C22").Value = '''38972 (56.6) '
$ws.Range("B23").Value = '''single'
$ws.Range("C23").Value = '''19326 (28.1) '
$ws.Range("B24").Value = '''widow/widower'
$ws.Range("C24").Value = '''10561 (15.3) '
$ws.Range("A25").Value = '''ECI_index_walraven (%)'
$ws.Range("B25").Value = '''-14'
$ws.Range("C25").Value = '''    1 ( 0.0) '
$ws.Range("B26").Value = '''-11'
$ws.Range("C26").Value = '''    1 ( 0.0) '
$ws.Range("B27").Value = '''-10'
$ws.Range("C27").Value = '''   10 ( 0.0) '
$ws.Range("B28").Value = '''-7'
$ws.Range("C28").Value = '''   38 ( 0.1) '
$ws.Range("B29").Value = '''-6'
$ws.Range("C29").Value = '''    5 ( 0.0) '
$ws.Range("B30").Value = '''-5'
$ws.Range("C30").Value = '''   15 ( 0.0) '
$ws.Range("B31").Value = '''-4'
$ws.Range("C31").Value = '''  416 ( 0.6) '
$ws.Range("B32").Value = '''-3'
$ws.Range("C32").Value = '''  342 ( 0.5) '
$ws.Range("B33").Value = '''-2'
$ws.Range("C33").Value = '''  123 ( 0.2) '
$ws.Range("B34").Value = '''-1'
$ws.Range("C34").Value = '''  406 ( 0.6) '
$ws.Range("B35").Value = '''0'
$ws.Range("C35").Value = '''61121 (88.8) '
$ws.Range("B36").Value = '''1'
$ws.Range("C36").Value = '''  107 ( 0.2) '
$ws.Range("B37").Value = '''2'
$ws.Range("C37").Value = '''  230 ( 0.3) '
$ws.Range("B38").Value = '''3'
$ws.Range("C38").Value = ''' 1148 ( 1.7) '
$ws.Range("B39").Value = '''4'
$ws.Range("C39").Value = ''' 1126 ( 1.6) '
$ws.Range("B40").Value = '''5'
$ws.Range("C40").Value = ''' 1769 ( 2.6) '
$ws.Range("B41").Value = '''6'
$ws.Range("C41").Value = '''  339 ( 0.5) '
$ws.Range("B42").Value = '''7'
$ws.Range("C42").Value = '''  370 ( 0.5) '
$ws.Range("B43").Value = '''8'
$ws.Range("C43").Value = '''  161 ( 0.2) '
$ws.Range("B44").Value = '''9'
$ws.Range("C44").Value = '''  239 ( 0.3) '
$ws.Range("B45").Value = '''10'
$ws.Range("C45").Value = '''  112 ( 0.2) '
$ws.Range("B46").Value = '''11'
$ws.Range("C46").Value = '''  164 ( 0.2) '
$ws.Range("B47").Value = '''12'
$ws.Range("C47").Value = '''  285 ( 0.4) '
$ws.Range("B48").Value = '''13'
$ws.Range("C48").Value = '''   30 ( 0.0) '
$ws.Range("B49").Value = '''14'
$ws.Range("C49").Value = '''   49 ( 0.1) '
$ws.Range("B50").Value = '''15'
$ws.Range("C50").Value = '''   46 ( 0.1) '
$ws.Range("B51").Value = '''16'
$ws.Range("C51").Value = '''   93 ( 0.1) '
$ws.Range("B52").Value = '''17'
$ws.Range("C52").Value = '''   41 ( 0.1) '
$ws.Range("B53").Value = '''18'
$ws.Range("C53").Value = '''   13 ( 0.0) '
$ws.Range("B54").Value = '''19'
$ws.Range("C54").Value = '''   15 ( 0.0) '
$ws.Range("B55").Value = '''20'
$ws.Range("C55").Value = '''   12 ( 0.0) '
$ws.Range("B56").Value = '''21'
$ws.Range("C56").Value = '''   16 ( 0.0) '
$ws.Range("B57").Value = '''22'
$ws.Range("C57").Value = '''    4 ( 0.0) '
$ws.Range("B58").Value = '''23'
$ws.Range("C58").Value = '''    2 ( 0.0) '
$ws.Range("B59").Value = '''24'
$ws.Range("C59").Value = '''    3 ( 0.0) '
$ws.Range("B60").Value = '''25'
$ws.Range("C60").Value = '''    1 ( 0.0) '
$ws.Range("B61").Value = '''26'
$ws.Range("C61").Value = '''    1 ( 0.0) '
$ws.Range("B62").Value = '''27'
$ws.Range("C62").Value = '''    2 ( 0.0) '
$ws.Range("B63").Value = '''28'
$ws.Range("C63").Value = '''    1 ( 0.0) '
$ws.Range("B64").Value = '''32'
$ws.Range("C64").Value = '''    1 ( 0.0) '
$ws.Range("B65").Value = '''33'
$ws.Range("C65").Value = '''    1 ( 0.0) '
$ws.Range("A66").Value = '''CCI_index_quan_original (%)'
$ws.Range("B66").Value = '''0'
$ws.Range("C66").Value = '''60970 (88.5) '
$ws.Range("B67").Value = '''1'
$ws.Range("C67").Value = ''' 4691 ( 6.8) '
$ws.Range("B68").Value = '''2'
$ws.Range("C68").Value = ''' 2161 ( 3.1) '
$ws.Range("B69").Value = '''3'
$ws.Range("C69").Value = '''  609 ( 0.9) '
$ws.Range("B70").Value = '''4'
$ws.Range("C70").Value = '''  195 ( 0.3) '
$ws.Range("B71").Value = '''5'
$ws.Range("C71").Value = '''   88 ( 0.1) '
$ws.Range("B72").Value = '''6'
$ws.Range("C72").Value = '''   40 ( 0.1) '
$ws.Range("B73").Value = '''7'
$ws.Range("C73").Value = '''   14 ( 0.0) '
$ws.Range("B74").Value = '''8'
$ws.Range("C74").Value = '''   60 ( 0.1) '
$ws.Range("B75").Value = '''9'
$ws.Range("C75").Value = '''   23 ( 0.0) '
$ws.Range("B76").Value = '''10'
$ws.Range("C76").Value = '''    2 ( 0.0) '
$ws.Range("B77").Value = '''11'
$ws.Range("C77").Value = '''    4 ( 0.0) '
$ws.Range("B78").Value = '''12'
$ws.Range("C78").Value = '''    2 ( 0.0) '
$ws.Range("A79").Value = '''Rx_index_index (%)'
$ws.Range("B79").Value = '''0'
$ws.Range("C79").Value = '''35849 (52.1) '
$ws.Range("B80").Value = '''1'
$ws.Range("C80").Value = ''' 4979 ( 7.2) '
$ws.Range("B81").Value = '''2'
$ws.Range("C81").Value = ''' 5437 ( 7.9) '
$ws.Range("B82").Value = '''3'
$ws.Range("C82").Value = ''' 5121 ( 7.4) '
$ws.Range("B83").Value = '''4'
$ws.Range("C83").Value = ''' 4718 ( 6.9) '
$ws.Range("B84").Value = '''5'
$ws.Range("C84").Value = ''' 4170 ( 6.1) '
$ws.Range("B85").Value = '''6'
$ws.Range("C85").Value = ''' 3297 ( 4.8) '
$ws.Range("B86").Value = '''7'
$ws.Range("C86").Value = ''' 2326 ( 3.4) '
$ws.Range("B87").Value = '''8'
$ws.Range("C87").Value = ''' 1422 ( 2.1) '
$ws.Range("B88").Value = '''9'
$ws.Range("C88").Value = '''  772 ( 1.1) '
$ws.Range("B89").Value = '''10'
$ws.Range("C89").Value = '''  420 ( 0.6) '
$ws.Range("B90").Value = '''11'
$ws.Range("C90").Value = '''  202 ( 0.3) '
$ws.Range("B91").Value = '''12'
$ws.Range("C91").Value = '''   85 ( 0.1) '
$ws.Range("B92").Value = '''13'
$ws.Range("C92").Value = '''   31 ( 0.0) '
$ws.Range("B93").Value = '''14'
$ws.Range("C93").Value = '''   21 ( 0.0) '
$ws.Range("B94").Value = '''15'
$ws.Range("C94").Value = '''    6 ( 0.0) '
$ws.Range("B95").Value = '''16'
$ws.Range("C95").Value = '''    2 ( 0.0) '
$ws.Range("B96").Value = '''17'
$ws.Range("C96").Value = '''    1 ( 0.0) '
